$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric but must stay text (matches original inlineStr type)
# Pre-format as Text so Excel does not coerce the literal string into a Number.
$textCells = @("D5", "D6", "D7", "D9", "D10", "D12", "D13", "D15", "D17", "D19", "D21", "D23", "D24", "D26", "D27", "D29", "D30", "D31", "D33", "D35", "D41", "D43", "D44", "D47", "D48", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply all cell value changes
$ws.Range("D2").Value = "52.049.03"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").Value = "3.004.27"
$ws.Range("E3").Value = "  +2.93%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "354.58"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").Value = "107.06"
$ws.Range("E6").Value = "  -3.02%  "
$ws.Range("D7").Value = "0.558"
$ws.Range("E7").Value = "  -1.52%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").Value = "0.612"
$ws.Range("E9").Value = "  -2.54%  "
$ws.Range("D10").Value = "38.13"
$ws.Range("E10").Value = "  -2.83%  "
$ws.Range("E11").Value = "  +2.36%  "
$ws.Range("D12").Value = "0.0857"
$ws.Range("E12").Value = "  -3.57%  "
$ws.Range("D13").Value = "19.03"
$ws.Range("E13").Value = "  -3.36%  "
$ws.Range("D14").Value = "3.478.10"
$ws.Range("E14").Value = "  +2.76%  "
$ws.Range("D15").Value = "7.65"
$ws.Range("E15").Value = "  -3.13%  "
$ws.Range("D16").Value = "2.997.40"
$ws.Range("E16").Value = "  +2.56%  "
$ws.Range("D17").Value = "1.02"
$ws.Range("E17").Value = "  +3.92%  "
$ws.Range("D18").Value = "52.095.98"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").Value = "3.40"
$ws.Range("E19").Value = "  +3.79%  "
$ws.Range("E20").Value = "  -1.40%  "
$ws.Range("D21").Value = "13.57"
$ws.Range("E21").Value = "  -3.47%  "
$ws.Range("D22").Value = "0.0₃0972"
$ws.Range("E22").Value = "  -0.96%  "
$ws.Range("D23").Value = "69.14"
$ws.Range("E23").Value = "  -2.40%  "
$ws.Range("D24").Value = "263.89"
$ws.Range("E24").Value = "  -2.13%  "
$ws.Range("E25").Value = "  -3.40%  "
$ws.Range("D26").Value = "0.179"
$ws.Range("E26").Value = "  -2.21%  "
$ws.Range("D27").Value = "27.04"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").Value = "7.37"
$ws.Range("E29").Value = "  -0.57%  "
$ws.Range("D30").Value = "0.107"
$ws.Range("E30").Value = "  -0.57%  "
$ws.Range("D31").Value = "6.39"
$ws.Range("E31").Value = "  +5.04%  "
$ws.Range("E32").Value = "  -3.69%  "
$ws.Range("D33").Value = "36.06"
$ws.Range("E33").Value = "  -5.56%  "
$ws.Range("E34").Value = "  +15.43%  "
$ws.Range("D35").Value = "51.07"
$ws.Range("E35").Value = "  -2.49%  "
$ws.Range("E36").Value = "  -0.67%  "
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("E38").Value = "  +2.82%  "
$ws.Range("E39").Value = "  +4.14%  "
$ws.Range("E40").Value = "  -2.61%  "
$ws.Range("D41").Value = "17.64"
$ws.Range("E41").Value = "  -4.22%  "
$ws.Range("E42").Value = "  -1.65%  "
$ws.Range("D43").Value = "124.58"
$ws.Range("E43").Value = "  +1.86%  "
$ws.Range("D44").Value = "22.80"
$ws.Range("E44").Value = "  -0.82%  "
$ws.Range("E45").Value = "  -1.92%  "
$ws.Range("D46").Value = "2.122.34"
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("D47").Value = "3.33"
$ws.Range("E47").Value = "  -3.25%  "
$ws.Range("D48").Value = "2.35"
$ws.Range("E48").Value = "  -6.45%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "3.304.49"
$ws.Range("E49").Value = "  +2.84%  "
$ws.Range("B50").Value = "TheGraph"
$ws.Range("C50").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D50").Value = "0.242"
$ws.Range("E50").Value = "  -3.76%  "
$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D51").Value = "0.0334"
$ws.Range("E51").Value = "  +0.42%  "

# Restore default style on the cells we temporarily reformatted, so only the
# cell content changes (no stray style attribute left on the cell).
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
